$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save current (pre-edit) values of row 13 and row 14 for the columns that change
# (use Value2 - the Value getter in this engine does not resolve properly)
$A13 = $ws.Range("A13").Value2
$B13 = $ws.Range("B13").Value2
$D13 = $ws.Range("D13").Value2
$E13 = $ws.Range("E13").Value2
$F13 = $ws.Range("F13").Value2
$G13 = $ws.Range("G13").Value2
$H13 = $ws.Range("H13").Value2
$Q13 = $ws.Range("Q13").Value2
$R13 = $ws.Range("R13").Value2

$A14 = $ws.Range("A14").Value2
$B14 = $ws.Range("B14").Value2
$D14 = $ws.Range("D14").Value2
$E14 = $ws.Range("E14").Value2
$F14 = $ws.Range("F14").Value2
$G14 = $ws.Range("G14").Value2
$H14 = $ws.Range("H14").Value2
$Q14 = $ws.Range("Q14").Value2
$R14 = $ws.Range("R14").Value2

# Write row 14's former values into row 13
$ws.Range("A13").Value2 = $A14
$ws.Range("B13").Value2 = $B14
$ws.Range("D13").Value2 = $D14
$ws.Range("E13").Value2 = $E14
$ws.Range("F13").Value2 = $F14
$ws.Range("G13").Value2 = $G14
$ws.Range("H13").Value2 = $H14
$ws.Range("Q13").Value2 = $Q14
$ws.Range("R13").Value2 = $R14

# Write row 13's former values into row 14
$ws.Range("A14").Value2 = $A13
$ws.Range("B14").Value2 = $B13
$ws.Range("D14").Value2 = $D13
$ws.Range("E14").Value2 = $E13
$ws.Range("F14").Value2 = $F13
$ws.Range("G14").Value2 = $G13
$ws.Range("H14").Value2 = $H13
$ws.Range("Q14").Value2 = $Q13
$ws.Range("R14").Value2 = $R13
